$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $replace, 2)
}

Replace-Text "Kontakt[[PHONE_2]], [[EMAIL_1]]" "Kontakt: +420 [[AMOUNT_1]], [[EMAIL_1]]"
Replace-Text "[[PHONE_3]]" "Tel.: [[AMOUNT_2]]"
Replace-Text "[[PHONE_4]]" "Telefon: +420 [[AMOUNT_3]]"
Replace-Text "částku 420 000 Kč." "částku [[AMOUNT_4]]."
Replace-Text "Kontakt: [[EMAIL_3]], [[PHONE_5]]" "Kontakt: [[EMAIL_3]], tel. [[AMOUNT_5]]"
Replace-Text "[[PHONE_6]]" "Mobil: +420 [[AMOUNT_6]]"
Replace-Text "nájemné 18 500 Kč." "nájemné [[AMOUNT_7]]."
Replace-Text "[[PHONE_7]]" "Tel.: +420 [[AMOUNT_8]]"
Replace-Text "[[PHONE_8]]" "Telefon: [[AMOUNT_9]]"
Replace-Text "[[PHONE_9]]" "Tel.: [[AMOUNT_10]]"
Replace-Text "Kontakt: [[EMAIL_6]],[[PHONE_10]]" "Kontakt: [[EMAIL_6]], [[AMOUNT_11]]"
Replace-Text "částku 32 000 Kč na účet" "částku [[AMOUNT_12]] na účet"
Replace-Text "[[PHONE_11]]" "Tel.: +420 [[AMOUNT_13]]"
Replace-Text "osobní kontakt: [[EMAIL_8]],[[PHONE_12]])" "osobní kontakt: [[EMAIL_8]], [[AMOUNT_14]])"
Replace-Text "kontakt[[PHONE_13]])." "kontakt: [[AMOUNT_15]])."
Replace-Text "[[PHONE_14]]" "Telefon: +420 [[AMOUNT_16]]"
Replace-Text "[[PHONE_15]]" "Tel.: [[AMOUNT_17]]"
Replace-Text "[[PHONE_16]]" "Telefon: [[AMOUNT_18]]"
Replace-Text "výši 220 000 Kč." "výši [[AMOUNT_19]]."
